$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $val
    $ws.Range($addr).Style = "Normal"
}

Set-TextValue "D2" "304.22"
Set-TextValue "E2" "0.63%"
Set-TextValue "D3" "36.01"
Set-TextValue "E3" "-3.84%"
Set-TextValue "D4" "5.097"
Set-TextValue "E4" "1.82%"
Set-TextValue "D5" "0.07843"
Set-TextValue "D6" "2.150"
Set-TextValue "E6" "-3.36%"
Set-TextValue "D7" "7.940"
Set-TextValue "E7" "-1.03%"
Set-TextValue "B8" "MXToken"
Set-TextValue "C8" "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue "D8" "0.9189"
Set-TextValue "E8" "1.16%"
Set-TextValue "B9" "LiechtensteinCryptoassetsExchange"
Set-TextValue "C9" "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextValue "D9" "0.09660"
Set-TextValue "E9" "0.39%"
Set-TextValue "B10" "WazirX"
Set-TextValue "C10" "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextValue "D10" "0.1863"
Set-TextValue "E10" "-1.65%"
Set-TextValue "B11" "MandalaExchangeToken"
Set-TextValue "C11" "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextValue "D11" "0.08584"
Set-TextValue "E11" "1.12%"
Set-TextValue "B12" "BitrueCoin"
Set-TextValue "C12" "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextValue "D12" "0.03501"
Set-TextValue "E12" "-0.74%"
Set-TextValue "B13" "BitMartToken"
Set-TextValue "C13" "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextValue "D13" "0.09935"
Set-TextValue "E13" "-0.32%"
Set-TextValue "B14" "BitForexToken"
Set-TextValue "C14" "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextValue "D14" "0.001435"
Set-TextValue "E14" "-3.79%"
Set-TextValue "B15" "TigerCash"
Set-TextValue "C15" "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextValue "D15" "0.005714"
Set-TextValue "E15" "1.26%"
Set-TextValue "B16" "LEO"
Set-TextValue "C16" "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue "D16" "3.463"
Set-TextValue "E16" "-0.07%"
Set-TextValue "B17" "GateToken"
Set-TextValue "C17" "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
Set-TextValue "D17" "4.104"
Set-TextValue "E17" "2.05%"
Set-TextValue "D18" "2.483"
Set-TextValue "E18" "20.18%"
Set-TextValue "D19" "0.3427"
Set-TextValue "E19" "-1.06%"
Set-TextValue "D20" "0.1320"
Set-TextValue "E20" "0.81%"
Set-TextValue "D21" "4.796"
Set-TextValue "E21" "0.71%"
Set-TextValue "D22" "0.2204"
Set-TextValue "E22" "-0.03%"
Set-TextValue "D23" "0.04546"
Set-TextValue "E23" "-2.24%"
Set-TextValue "D24" "0.005083"
Set-TextValue "E24" "14.20%"
Set-TextValue "D25" "0.001240"
Set-TextValue "E25" "0.87%"
Set-TextValue "D27" "0.0004759"
Set-TextValue "E27" "0.12%"
Set-TextValue "D39" "0.01842"
Set-TextValue "E39" "4.94%"
Set-TextValue "D40" "0.04723"
Set-TextValue "E40" "0.04%"
Set-TextValue "D41" "0.007739"
Set-TextValue "E41" "-1.49%"
Set-TextValue "D42" "0.1401"
Set-TextValue "E42" "0.64%"
Set-TextValue "D43" "0.007742"
Set-TextValue "E43" "1.01%"
Set-TextValue "D44" "0.002242"
Set-TextValue "E44" "3.29%"
Set-TextValue "D45" "0.01120"
Set-TextValue "E45" "13.33%"
Set-TextValue "D46" "0.00006417"
Set-TextValue "E46" "5.45%"
Set-TextValue "D47" "0.00000000752"
Set-TextValue "D48" "0.0005812"
Set-TextValue "E48" "0.19%"
Set-TextValue "D49" "36.35"
Set-TextValue "E49" "319.25%"
Set-TextValue "D50" "0.002004"
Set-TextValue "E50" "-25.55%"
Set-TextValue "D51" "0.00002104"
